# Updated symbol list on Tue Dec 13 21:57:53 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as text (numeric-looking strings),
# so force the cells to a text number format before assigning, otherwise
# Excel would auto-coerce the numeric-looking text into a real number.
$priceCells = @(
    "D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12",
    "D14","D16","D18","D19","D20","D21","D22","D23","D24",
    "D40","D41","D42","D43","D45","D48","D49","D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value  = "270.56"
$ws.Range("D3").Value  = "23.04"
$ws.Range("D4").Value  = "6.321"
$ws.Range("D5").Value  = "0.06202"
$ws.Range("D6").Value  = "3.639"
$ws.Range("D7").Value  = "6.685"
$ws.Range("D8").Value  = "1.391"
$ws.Range("D9").Value  = "0.8295"
$ws.Range("D10").Value = "0.01379"
$ws.Range("D11").Value = "0.1600"
$ws.Range("D12").Value = "0.08285"
$ws.Range("D14").Value = "0.03192"
$ws.Range("D16").Value = "3.867"
$ws.Range("D18").Value = "0.04734"
$ws.Range("D19").Value = "0.006364"
$ws.Range("D20").Value = "0.005680"
$ws.Range("D21").Value = "0.001077"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D23").Value = "3.717"
$ws.Range("D24").Value = "2.413"
$ws.Range("D40").Value = "0.04696"
$ws.Range("D41").Value = "0.007024"
$ws.Range("D45").Value = "0.00006258"
$ws.Range("D48").Value = "0.9196"
$ws.Range("D49").Value = "0.001377"
$ws.Range("D50").Value = "0.00001399"

# Rows 42 and 43 swap their coin identity (CEJI <-> BKEXToken), with
# refreshed price/volume data.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1162"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003299"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Write-Host "Applied cryptos.xlsx price/symbol updates"
